$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report date range) ---
$ws.Range("A8").Value = "Volume 30   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Range("N14").Value = -71.428571428571
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 400
$ws.Range("L15").Value = 400
$ws.Range("M15").Value = 66.666666666666
$ws.Range("N15").Value = -44.444444444444
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 32
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = -5.882352941176
$ws.Range("M16").Value = -41.818181818181
$ws.Range("N16").Value = -84.615384615384
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 16
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 39
$ws.Range("J17").Value = 35
$ws.Range("K17").Value = 11.428571428571
$ws.Range("L17").Value = -13.333333333333
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = -61.386138613861
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 6.666666666666
$ws.Range("I18").Value = 40
$ws.Range("J18").Value = 68
$ws.Range("K18").Value = -41.176470588235
$ws.Range("L18").Value = 5.263157894736
$ws.Range("M18").Value = 17.647058823529
$ws.Range("N18").Value = -87.220447284345
$ws.Range("C19").Value = 3
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -76.923076923076
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 5.555555555555
$ws.Range("I19").Value = 140
$ws.Range("J19").Value = 128
$ws.Range("K19").Value = 9.375
$ws.Range("L19").Value = 44.329896907216
$ws.Range("M19").Value = 7.692307692307
$ws.Range("N19").Value = -49.458483754512
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 40
$ws.Range("I20").Value = 29
$ws.Range("K20").Value = 93.333333333333
$ws.Range("L20").Value = 123.076923076923
$ws.Range("M20").Value = 190
$ws.Range("N20").Value = -88.306451612903
$ws.Range("C21").Value = 13
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -38.095238095238
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = 23.611111111111
$ws.Range("I21").Value = 287
$ws.Range("J21").Value = 287
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 25.877192982456
$ws.Range("M21").Value = 9.125475285171
$ws.Range("N21").Value = -75.322441960447
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 5
$ws.Range("J22").Value = 9
$ws.Range("K22").Value = -44.444444444444
$ws.Range("L22").Value = 25
$ws.Range("M22").Value = -44.444444444444
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -25
$ws.Range("F23").Value = 12
$ws.Range("G23").Value = 8
$ws.Range("I23").Value = 30
$ws.Range("J23").Value = 28
$ws.Range("K23").Value = 7.142857142857
$ws.Range("L23").Value = -26.829268292682
$ws.Range("M23").Value = 7.142857142857
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 35
$ws.Range("E24").Value = -37.142857142857
$ws.Range("F24").Value = 106
$ws.Range("G24").Value = 170
$ws.Range("H24").Value = -37.647058823529
$ws.Range("I24").Value = 493
$ws.Range("J24").Value = 532
$ws.Range("K24").Value = -7.330827067669
$ws.Range("L24").Value = 108.016877637131
$ws.Range("M24").Value = 92.578125
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 12.5
$ws.Range("F25").Value = 23
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = 4.545454545454
$ws.Range("I25").Value = 73
$ws.Range("J25").Value = 70
$ws.Range("K25").Value = 4.285714285714
$ws.Range("L25").Value = 15.873015873015
$ws.Range("M25").Value = -7.594936708860
$ws.Range("F26").Value = 2
$ws.Range("I26").Value = 6
$ws.Range("K26").Value = 200
$ws.Range("L26").Value = 100
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 20
$ws.Range("J27").Value = 13
$ws.Range("K27").Value = 30.769230769230
$ws.Range("L27").Value = 30.769230769230
$ws.Range("N28").Value = -90.909090909090
$ws.Range("N29").Value = -90.909090909090
$ws.Range("L30").Value = -80

# --- Cells changing between numeric and text representation ---
# Strategy: set the new value first (so the underlying type is correct;
# a leading apostrophe forces text interpretation for numeric-looking
# strings like "0"), then copy *formatting only* from a stable donor cell
# that already has the desired style index, to avoid growing styles.xml
# with ad-hoc formats.

$ws.Range("C15").Value = 1
$ws.Range("D25").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("C16").Value = 2
$ws.Range("D25").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("D20").Value = "'0"
$ws.Range("C28").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").Value = "***.*"
$ws.Range("C28").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("C22").Value = 1
$ws.Range("D25").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D22").Value = 1
$ws.Range("D25").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = 0
$ws.Range("K23").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("C26").Value = 1
$ws.Range("D25").Copy()
$ws.Range("C26").PasteSpecial(-4122)

$ws.Range("C27").Value = "'0"
$ws.Range("C28").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").Value = 1
$ws.Range("D25").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = -100
$ws.Range("K23").Copy()
$ws.Range("E27").PasteSpecial(-4122)
